# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-55) previously listed periods in
# descending order (2003 down to 1612). It is updated to list the periods
# in ascending / chronological order (1612 up to 2003), and the
# "Valor Mora" (F) / "Salario Basico" (G) figures for that range are
# refreshed to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chronological (ascending) list of period codes that now populate
# E16:E55, replacing the previous descending list.
$periods = @(
    "1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("E$row").Value = $periods[$i]

    if ($row -le 36) {
        $ws.Range("F$row").Value = 27578
    } else {
        $ws.Range("F$row").Value = 31249
    }

    $ws.Range("G$row").Value = 781242
}
